# Update "想去人数" (want-to-go count) figures in column F across the
# workbook's four sheets, matching the refreshed data snapshot.

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 1262
$ws1.Range("F4").Value  = 58
$ws1.Range("F5").Value  = 5579
$ws1.Range("F6").Value  = 1786
$ws1.Range("F7").Value  = 1786
$ws1.Range("F8").Value  = 6360
$ws1.Range("F12").Value = 8
$ws1.Range("F14").Value = 32
$ws1.Range("F18").Value = 7922
$ws1.Range("F19").Value = 7922
$ws1.Range("F29").Value = 47
$ws1.Range("F30").Value = 173
$ws1.Range("F31").Value = 1751
$ws1.Range("F38").Value = 3920

# Sheet 2: 演出 (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F10").Value = 2
$ws2.Range("F14").Value = 27

# Sheet 3: 本地生活 (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 2278

# Sheet 4: 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 2278
$ws4.Range("F5").Value  = 1262
$ws4.Range("F7").Value  = 58
$ws4.Range("F10").Value = 5579
$ws4.Range("F12").Value = 1786
$ws4.Range("F13").Value = 1786
$ws4.Range("F14").Value = 6360
$ws4.Range("F17").Value = 2
$ws4.Range("F21").Value = 32
$ws4.Range("F24").Value = 7922
$ws4.Range("F25").Value = 7922
$ws4.Range("F34").Value = 47
$ws4.Range("F35").Value = 173
$ws4.Range("F36").Value = 1751
$ws4.Range("F41").Value = 27
$ws4.Range("F47").Value = 3920
